# feat: add 2022-Q1 data
#
# The workbook previously tracked a running "总计" (totals) sheet as the
# 4th/last tab. This edit:
#   1. Turns the old "总计" tab into the per-fund holdings detail sheet for
#      the new quarter, renamed "2022-Q1", with its own header row/data
#      (same layout as the 2020-Q4 / 2021-Q1 / 2021-Q2 detail sheets).
#   2. Appends a brand-new "总计" tab at the end of the workbook containing
#      the historical summary rows plus a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# Existing detail sheet whose header / row-index styling (bold, thin
# border, centered) matches what the "总计" sheet used, and which stays
# untouched by this edit - used below as a format-paste source so the
# repurposed/new sheets reuse the exact same styles already in this
# workbook instead of ad-hoc COM-built ones.
$styleSource = $wb.Worksheets.Item("2021-Q2")

# ---------------------------------------------------------------------
# Step 1: repurpose the existing "总计" sheet as the "2022-Q1" detail sheet
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A1:D4").Clear()
$totalSheet.Name = "2022-Q1"

$totalSheet.Range("B1").Value = "基金代码"
$totalSheet.Range("C1").Value = "基金名称"
$totalSheet.Range("D1").Value = "基金规模"
$totalSheet.Range("E1").Value = "股票总仓位"
$totalSheet.Range("F1").Value = "仓位占比"
$totalSheet.Range("G1").Value = "持有市值(亿元)"
$totalSheet.Range("H1").Value = "仓位排名"

# Apply the matching header style (copy formats only, values untouched)
$styleSource.Range("B1").Copy()
$totalSheet.Range("B1:H1").PasteSpecial(-4122)

# Row-index column (A) shares the bold/bordered/centered header look
$styleSource.Range("A2").Copy()
$totalSheet.Range("A2:A5").PasteSpecial(-4122)

# Numeric-looking text columns (D/E/F/G) must stay text, like "38.60"
$totalSheet.Range("D2:G5").NumberFormat = "@"
# Fund-code column (B) must stay text too, to keep leading zeros (014591)
$totalSheet.Range("B2:B5").NumberFormat = "@"

$fundRows = @(
    @{ A=0; B="014591"; C="广发瑞誉一年持有期混合A";             D="38.60"; E="93.40"; F="5.30"; G="2.0458"; H=4 },
    @{ A=1; B="014592"; C="广发瑞誉一年持有期混合C";             D="4.08";  E="93.40"; F="5.30"; G="0.2162"; H=4 },
    @{ A=2; B="004316"; C="前海开源沪港深裕鑫灵活配置混合A"; D="0.64";  E="90.55"; F="3.15"; G="0.0202"; H=4 },
    @{ A=3; B="004317"; C="前海开源沪港深裕鑫灵活配置混合C"; D="0.47";  E="90.55"; F="3.15"; G="0.0148"; H=4 }
)

$r = 2
foreach ($row in $fundRows) {
    $totalSheet.Cells.Item($r, 1).Value = $row.A
    $totalSheet.Cells.Item($r, 2).Value = $row.B
    $totalSheet.Cells.Item($r, 3).Value = $row.C
    $totalSheet.Cells.Item($r, 4).Value = $row.D
    $totalSheet.Cells.Item($r, 5).Value = $row.E
    $totalSheet.Cells.Item($r, 6).Value = $row.F
    $totalSheet.Cells.Item($r, 7).Value = $row.G
    $totalSheet.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: add a fresh "总计" sheet at the end with the historical summary
#         rows plus the new 2022-Q1 row on top
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$summarySheet = $wb.Worksheets.Add($null, $lastSheet)
$summarySheet.Name = "总计"

$summarySheet.Range("B1").Value = "日期"
$summarySheet.Range("C1").Value = "持有数量(只)"
$summarySheet.Range("D1").Value = "持有市值(亿元)"

# Apply the matching header style (copy formats only, values untouched)
$styleSource.Range("B1").Copy()
$summarySheet.Range("B1:D1").PasteSpecial(-4122)

# Row-index column (A) shares the bold/bordered/centered header look
$styleSource.Range("A2").Copy()
$summarySheet.Range("A2:A5").PasteSpecial(-4122)

$summaryRows = @(
    @{ A=0; B="2022-Q1"; C=4; D=2.3 },
    @{ A=1; B="2021-Q2"; C=5; D=6.16 },
    @{ A=2; B="2021-Q1"; C=2; D=1.85 },
    @{ A=3; B="2020-Q4"; C=4; D=2.07 }
)

$r = 2
foreach ($row in $summaryRows) {
    $summarySheet.Cells.Item($r, 1).Value = $row.A
    $summarySheet.Cells.Item($r, 2).Value = $row.B
    $summarySheet.Cells.Item($r, 3).Value = $row.C
    $summarySheet.Cells.Item($r, 4).Value = $row.D
    $r = $r + 1
}
